$d = $word.ActiveDocument

# Find the paragraph that contains the " m:'doc.html'.fromHTMLURI() " field
# (the one built from fldChar begin/end + instrText runs).
$fieldPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $fieldPara = $p
        break
    }
}

if ($fieldPara -eq $null) {
    throw "Could not locate the paragraph containing the M2Doc field."
}

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Rebuild the paragraph as plain literal text runs - "{", "m", ":", "'",
# "doc.html", <the _GoBack bookmark>, "'.fromHTMLURI()", "}" - instead of a
# real Word field (fldChar begin/end wrapping instrText runs). This is what
# TokenIteratorFieldRewriterSplit expects to parse as an M2Doc token.
$newParaXml =
    "<w:p $wNs w:rsidR='00C52979' w:rsidRDefault='00C52979' w:rsidP='00F5495F'>" +
        "<w:r><w:t>{</w:t></w:r>" +
        "<w:r><w:t>m</w:t></w:r>" +
        "<w:r><w:t>:</w:t></w:r>" +
        "<w:r><w:t>'</w:t></w:r>" +
        "<w:r><w:t>doc.html</w:t></w:r>" +
        "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
        "<w:bookmarkEnd w:id='0'/>" +
        "<w:r><w:t>'.fromHTMLURI()</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>}</w:t></w:r>" +
    "</w:p>"

# InsertXML replaces the full contents of the target range (including its
# own enclosing <w:p>), so calling it on the field paragraph's Range swaps
# the old field markup for the new literal-text run sequence.
$fieldPara.Range.InsertXML($newParaXml)
